$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped
# uniformly from 45192 (2023-09-23) to 45202 (2023-10-03) for every
# data row (rows 2 through 469).
$ws.Range("C2:C469").Value = 45202
